$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before K, shifting K:N -> L:O (mirrors the author's
# "commit báo cáo đơn giá" edit that adds a "Ngày thông báo" column).
$ws.Columns("K:K").Insert()

# Fill in the header for the newly inserted column (row 6 of the 2-row
# header block) with the new shared string.
$ws.Range("K6").Value = "Ngày thông báo"

# Match the author's saved selection.
$ws.Range("K7").Select()
